$wb = $excel.ActiveWorkbook

# Sheet: 展览
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 449
$wsExhibit.Range("F6").Value = 1297
$wsExhibit.Range("F7").Value = 500
$wsExhibit.Range("F9").Value = 274
$wsExhibit.Range("F12").Value = 1086
$wsExhibit.Range("F16").Value = 81
$wsExhibit.Range("F18").Value = 1636
$wsExhibit.Range("F19").Value = 603
$wsExhibit.Range("F20").Value = 265
$wsExhibit.Range("F21").Value = 173
$wsExhibit.Range("F22").Value = 1631
$wsExhibit.Range("F23").Value = 395
$wsExhibit.Range("F28").Value = 1902
$wsExhibit.Range("F29").Value = 2800
$wsExhibit.Range("F30").Value = 1583
$wsExhibit.Range("F33").Value = 620
$wsExhibit.Range("F35").Value = 1718
$wsExhibit.Range("F36").Value = 878
$wsExhibit.Range("F37").Value = 1749
$wsExhibit.Range("F38").Value = 192
$wsExhibit.Range("F40").Value = 828
$wsExhibit.Range("F42").Value = 824
$wsExhibit.Range("F43").Value = 778
$wsExhibit.Range("F44").Value = 982

# Sheet: 演出
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F12").Value = 784
$wsShow.Range("F13").Value = 24

# Sheet: 全部类型
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 449
$wsAll.Range("F7").Value = 1297
$wsAll.Range("F8").Value = 500
$wsAll.Range("F9").Value = 274
$wsAll.Range("F13").Value = 1086
$wsAll.Range("F18").Value = 1636
$wsAll.Range("F19").Value = 603
$wsAll.Range("F20").Value = 173
$wsAll.Range("F21").Value = 1631
$wsAll.Range("F23").Value = 395
$wsAll.Range("F28").Value = 2800
$wsAll.Range("F29").Value = 1583
$wsAll.Range("F32").Value = 784
$wsAll.Range("F33").Value = 24
$wsAll.Range("F36").Value = 620
$wsAll.Range("F37").Value = 1718
$wsAll.Range("F40").Value = 878
$wsAll.Range("F41").Value = 1749
$wsAll.Range("F42").Value = 828
$wsAll.Range("F43").Value = 824
$wsAll.Range("F44").Value = 778
$wsAll.Range("F45").Value = 982

